{"js": "// Load all paragraphs in the document body.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// 1) Remove the \"Meta description: ...\" paragraph that follows the title\n//    (Heading1) paragraph. It is identified by its leading bold\n//    \"Meta description\" run.\nlet metaPara = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  p.load(\"text\");\n}\nawait context.sync();\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text.indexOf(\"Meta description:\") === 0) {\n    metaPara = p;\n    break;\n  }\n}\nif (metaPara) {\n  metaPara.delete();\n  await context.sync();\n}\n\n// 2) Re-load the paragraph collection (indexes shifted after the delete)\n//    and locate the closing \"Please create a cartoon...\" image-prompt\n//    paragraph, which is the very last paragraph in the document body.\nconst paragraphs2 = context.document.body.paragraphs;\nparagraphs2.load(\"items\");\nawait context.sync();\n\nconst lastPara = paragraphs2.items[paragraphs2.items.length - 1];\n\n// 3) Insert a new paragraph right before it containing the page title,\n//    bold (matching the title used at the top of the document), and make\n//    sure it does not inherit the italic formatting of the neighboring\n//    paragraph.\nconst titlePara = lastPara.insertParagraph(\n  \"Play Da Vinci Ways for Free - Slot Game Review\",\n  Word.InsertLocation.before\n);\ntitlePara.font.set({ bold: true, italic: false });\nawait context.sync();\n\n// 4) Replace the text of the final (image-prompt) paragraph with the old\n//    meta-description body text, keeping its existing (italic)\n//    formatting intact.\nlastPara.insertText(\n  \"Experience the beauty of Da Vinci's art while potentially winning big in this slot game with Tumbling Reels and Free Spins. Play for free now.\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Remove the \"Meta description: ...\" paragraph that follows the title\n#    (Heading1) paragraph. Find it by its distinctive leading text rather\n#    than assuming a fixed index.\n$metaIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs($i).Range.Text\n    if ($t -like \"Meta description:*\") {\n        $metaIndex = $i\n        break\n    }\n}\nif ($metaIndex -gt 0) {\n    $d.Paragraphs($metaIndex).Range.Delete()\n}\n\n# 2) The closing \"Please create a cartoon...\" image-prompt paragraph is\n#    the very last paragraph in the document body. Insert a new\n#    paragraph right before it containing the page title, bold (matching\n#    the title used at the top of the document).\n$count = $d.Paragraphs.Count\n$lastPara = $d.Paragraphs($count)\n$lastPara.Range.InsertParagraphBefore()\n\n$titleIndex = $count\n$titlePara = $d.Paragraphs($titleIndex)\n$titleRange = $titlePara.Range\n$titleRange.MoveEnd(1, -1) | Out-Null\n$titleRange.Text = \"Play Da Vinci Ways for Free - Slot Game Review\"\n$titleRange2 = $titlePara.Range\n$titleRange2.MoveEnd(1, -1) | Out-Null\n$titleRange2.Font.Bold = 1\n$titleRange2.Font.Italic = 0\n\n# 3) Replace the text of the final (image-prompt) paragraph with the old\n#    meta-description body text, keeping its existing (italic)\n#    formatting intact. Shrink the range by one character first so the\n#    trailing paragraph mark (and its formatting) is left untouched.\n$finalIndex = $d.Paragraphs.Count\n$finalPara = $d.Paragraphs($finalIndex)\n$finalRange = $finalPara.Range\n$finalRange.MoveEnd(1, -1) | Out-Null\n$finalRange.Text = \"Experience the beauty of Da Vinci's art while potentially winning big in this slot game with Tumbling Reels and Free Spins. Play for free now.\"\n"}
